$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "24/10/2025"
$ws.Range("B7").Value = "Neom FC"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "Al Khaleej"
$ws.Range("F7").Value = "D"
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0.87
$ws.Range("L7").Value = 2.04
$ws.Range("M7").Value = 13
$ws.Range("N7").Value = 15
$ws.Range("O7").Value = 5
$ws.Range("P7").Value = 4
